# Weekly fruit/hortaliza update: a new weekly price observation is inserted
# as row 712 (Choclo, Dulce o Americano, Primera, Región de Arica y
# Parinacota), pushing the existing rows 712-759 down to 713-760.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 712, shifting rows 712..759 down to 713..760.
$ws.Rows.Item(712).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(712, 1).Value = 3
$ws.Cells.Item(712, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(712, 3).Value = "Coquimbo"
$ws.Cells.Item(712, 4).Value = 44783
$ws.Cells.Item(712, 5).Value = 5
$ws.Cells.Item(712, 6).Value = 100112024
$ws.Cells.Item(712, 7).Value = "Choclo"
$ws.Cells.Item(712, 8).Value = "Dulce o Americano"
$ws.Cells.Item(712, 9).Value = "Primera"
$ws.Cells.Item(712, 10).Value = 115
$ws.Cells.Item(712, 11).Value = 34000
$ws.Cells.Item(712, 12).Value = 35000
$ws.Cells.Item(712, 13).Value = 34478
$ws.Cells.Item(712, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(712, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(712, 16).Value = 493
$ws.Cells.Item(712, 17).Value = 70
$ws.Cells.Item(712, 18).Value = "Hortaliza"

# Match the date number format already used by the rest of column D.
$ws.Cells.Item(712, 4).NumberFormat = $ws.Cells.Item(711, 4).NumberFormat
